$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power Storage")

# Update ExisUnits (column E) values for rows 7, 9, 10
$ws.Range("E7").Value = 5
$ws.Range("E9").Value = 2
$ws.Range("E10").Value = 29

# Update MaxInvest (column S) values for rows 7-11
$ws.Range("S7").Value = 8
$ws.Range("S8").Value = 8
$ws.Range("S9").Value = 8
$ws.Range("S10").Value = 8
$ws.Range("S11").Value = 8

# Update the selection in the frozen (bottomLeft) pane to S12
$ws.Activate()
$ws.Range("S12").Select()
